# Adds two new input tabs used by the `to_report` function:
#   - case_text_elements     (between "configurations" and "key_outputs")
#   - generic_text_elements  (appended as the last tab)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) case_text_elements — inserted right after "configurations"
# ---------------------------------------------------------------------------
$configSheet = $wb.Worksheets.Item("configurations")
$caseSheet = $wb.Worksheets.Add($null, $configSheet)
$caseSheet.Name = "case_text_elements"

# Column A first, then column B (matches how the data was actually typed in)
$caseSheet.Range("A1").Value = "case_text_element"
$caseSheet.Range("A2").Value = "title_front_page"
$caseSheet.Range("A3").Value = "strategic_challenge"

$caseSheet.Range("B1").Value = "value"
$caseSheet.Range("B2").Value = "Report of the DSM case"
$caseSheet.Range("B3").Value = "How to source energy?"

$caseHeader = $caseSheet.Range("A1:B1")
$caseHeader.Font.Bold = $true
$caseHeader.Borders.LineStyle = 1
$caseHeader.HorizontalAlignment = -4108
$caseHeader.VerticalAlignment = -4160

$caseSheet.Columns.Item(1).ColumnWidth = 19.7265625
$caseSheet.Columns.Item(2).ColumnWidth = 23.08984375

$caseSheet.Activate()

# ---------------------------------------------------------------------------
# 2) generic_text_elements — appended as the final tab
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$genSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$genSheet.Name = "generic_text_elements"

$genData = @(
  @("generic_text_element", "value"),
  @("title_strategic_challenge", "Strategic Challenge"),
  @("title_key_outputs", "Key outputs"),
  @("title_dmo", "Decision makers options (DMOs)"),
  @("title_scenarios", "Scenarios"),
  @("title_fixed_inputs", "Fixed inputs"),
  @("title_dependency_graph", "Dependency graph"),
  @("title_weighted_graph", "Resulting appreciations of different DMOs for scenario: "),
  @("intro_key_outputs", "The outputs upon which the decision makers will base their decision. Key outputs are often referred to as KPIs. Key outputs are grouped into themes."),
  @("intro_dmo", "Decision makers options are formulated by assigning a single value to all internal variable inputs. These inputs can be formulated and determined by the decision makers."),
  @("intro_scenarios", "Each external variable input can be thought of as a single aspect of external uncertainty affecting the outcome of the decision in scope. A scenario is defined by assigning a single value to all external variable inputs."),
  @("intro_fixed_inputs", "The inputs which only takes one value for all scenarios."),
  @("intro_dependency_graph", ""),
  @("intro_weighted_graph", "")
)

for ($i = 0; $i -lt $genData.Count; $i++) {
  $rowNum = $i + 1
  $key = $genData[$i][0]
  $val = $genData[$i][1]
  $genSheet.Cells.Item($rowNum, 1).Value = $key
  if ($val -ne "") {
    $genSheet.Cells.Item($rowNum, 2).Value = $val
  }
}

$genFullRange = $genSheet.Range("A1:B14")
$genFullRange.Font.Name = "Calibri"
$genFullRange.Font.Size = 11

$genHeader = $genSheet.Range("A1:B1")
$genHeader.Font.Bold = $true
$genHeader.Borders.LineStyle = 1
$genHeader.HorizontalAlignment = -4108
$genHeader.VerticalAlignment = -4160

$genSheet.Columns.Item(1).ColumnWidth = 23.90625
$genSheet.Columns.Item(2).ColumnWidth = 38.08984375

$genSheet.Range("B8").Select()

# Restore the originally active tab as in the published workbook
$caseSheet.Activate()
